# Auto-generated Excel COM-interop script applying scheduled market-price refresh
# to the leve profit columns (H..N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 7168.3335
$ws.Cells.Item(98, 10).Value = 20000
$ws.Cells.Item(98, 12).Value = 20000
$ws.Cells.Item(98, 14).Value = -22996

$ws.Cells.Item(122, 8).Value = 7168.3335
$ws.Cells.Item(122, 10).Value = 20000
$ws.Cells.Item(122, 12).Value = 60000
$ws.Cells.Item(122, 14).Value = -64900

$ws.Cells.Item(132, 8).Value = 1138.1
$ws.Cells.Item(132, 9).Value = 1138.1
$ws.Cells.Item(132, 11).Value = 3414.3
$ws.Cells.Item(132, 13).Value = -884.2999999999997

$ws.Cells.Item(138, 8).Value = 3644.4375
$ws.Cells.Item(138, 10).Value = 4434.5835
$ws.Cells.Item(138, 12).Value = 13303.7505
$ws.Cells.Item(138, 14).Value = -23583.7505

$ws.Cells.Item(141, 8).Value = 1000.75
$ws.Cells.Item(141, 10).Value = 999
$ws.Cells.Item(141, 12).Value = 2997
$ws.Cells.Item(141, 14).Value = -13357

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 12).ClearContents()
$ws.Cells.Item(4, 13).ClearContents()

$ws.Cells.Item(5, 8).Value = 250
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 250
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).Value = 250
$ws.Cells.Item(5, 14).Value = -474
$ws.Cells.Item(5, 12).ClearContents()

$ws.Cells.Item(12, 8).Value = 49999.5
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 49999.5
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 13).Value = 49999.5
$ws.Cells.Item(12, 14).Value = -50345.5
$ws.Cells.Item(12, 12).ClearContents()

$ws.Cells.Item(74, 8).Value = 21417.084
$ws.Cells.Item(74, 9).Value = 20609.174
$ws.Cells.Item(74, 11).Value = 20609.174
$ws.Cells.Item(74, 13).Value = -19735.174

$ws.Cells.Item(77, 8).Value = 21417.084
$ws.Cells.Item(77, 9).Value = 20609.174
$ws.Cells.Item(77, 11).Value = 103045.87
$ws.Cells.Item(77, 13).Value = -98677.87

$ws.Cells.Item(122, 8).Value = 1985.6666
$ws.Cells.Item(122, 9).Value = 1985.6666
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 5956.9998
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).Value = -3506.9998
$ws.Cells.Item(122, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 2261.682
$ws.Cells.Item(132, 9).Value = 1810
$ws.Cells.Item(132, 11).Value = 5430
$ws.Cells.Item(132, 13).Value = -2900

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 250
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 250
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).Value = 250
$ws.Cells.Item(4, 14).Value = -480
$ws.Cells.Item(4, 12).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1772.4286
$ws.Cells.Item(16, 9).Value = 1734.5
$ws.Cells.Item(16, 11).Value = 1734.5
$ws.Cells.Item(16, 13).Value = -1447.5

$ws.Cells.Item(31, 8).Value = 2499.6667
$ws.Cells.Item(31, 9).Value = 2600
$ws.Cells.Item(31, 11).Value = 2600
$ws.Cells.Item(31, 13).Value = -2305

$ws.Cells.Item(34, 8).Value = 2499.6667
$ws.Cells.Item(34, 9).Value = 2600
$ws.Cells.Item(34, 11).Value = 2600
$ws.Cells.Item(34, 13).Value = -2398

$ws.Cells.Item(113, 8).Value = 1772.4286
$ws.Cells.Item(113, 9).Value = 1734.5
$ws.Cells.Item(113, 11).Value = 1734.5
$ws.Cells.Item(113, 13).Value = 435.5

$ws.Cells.Item(132, 8).Value = 1778.875
$ws.Cells.Item(132, 10).Value = 2706.1667
$ws.Cells.Item(132, 12).Value = 8118.500100000001
$ws.Cells.Item(132, 14).Value = -13178.5001

$ws.Cells.Item(134, 8).Value = 5154.077
$ws.Cells.Item(134, 9).Value = 4521.5
$ws.Cells.Item(134, 11).Value = 13564.5
$ws.Cells.Item(134, 13).Value = -11029.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 370.1
$ws.Cells.Item(12, 9).Value = 478.2
$ws.Cells.Item(12, 10).Value = 262
$ws.Cells.Item(12, 11).Value = 1434.6
$ws.Cells.Item(12, 12).Value = 786
$ws.Cells.Item(12, 13).Value = -1261.6
$ws.Cells.Item(12, 14).Value = -1132

$ws.Cells.Item(107, 8).Value = 573
$ws.Cells.Item(107, 10).Value = 601.6
$ws.Cells.Item(107, 12).Value = 1804.8
$ws.Cells.Item(107, 14).Value = -5644.8

$ws.Cells.Item(132, 8).Value = 12547.5
$ws.Cells.Item(132, 10).Value = 45000
$ws.Cells.Item(132, 12).Value = 405000
$ws.Cells.Item(132, 14).Value = -410060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 14).Value = 0
$ws.Cells.Item(70, 12).ClearContents()

$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 14).Value = 0
$ws.Cells.Item(73, 12).ClearContents()

$ws.Cells.Item(113, 8).Value = 5412.5
$ws.Cells.Item(113, 9).Value = 5412.5
$ws.Cells.Item(113, 11).Value = 5412.5
$ws.Cells.Item(113, 13).Value = -3242.5

$ws.Cells.Item(132, 8).Value = 1641.2632
$ws.Cells.Item(132, 9).Value = 1290.75
$ws.Cells.Item(132, 11).Value = 3872.25
$ws.Cells.Item(132, 13).Value = -1342.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1393.6666
$ws.Cells.Item(55, 10).Value = 462.5
$ws.Cells.Item(55, 12).Value = 462.5
$ws.Cells.Item(55, 14).Value = -808.5

$ws.Cells.Item(93, 8).Value = 2825.75
$ws.Cells.Item(93, 9).Value = 2601
$ws.Cells.Item(93, 10).Value = 3500
$ws.Cells.Item(93, 11).Value = 2601
$ws.Cells.Item(93, 12).Value = 3500
$ws.Cells.Item(93, 13).Value = -1353
$ws.Cells.Item(93, 14).Value = -5996

$ws.Cells.Item(100, 8).Value = 4750
$ws.Cells.Item(100, 9).Value = 4750
$ws.Cells.Item(100, 11).Value = 4750
$ws.Cells.Item(100, 13).Value = -4209

$ws.Cells.Item(119, 8).Value = 78210
$ws.Cells.Item(119, 10).Value = 78210
$ws.Cells.Item(119, 12).Value = 78210
$ws.Cells.Item(119, 14).Value = -87886

$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 14).Value = 0
$ws.Cells.Item(121, 12).ClearContents()

$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 14).Value = 0
$ws.Cells.Item(122, 12).ClearContents()
$ws.Cells.Item(122, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 2928.7827
$ws.Cells.Item(132, 9).Value = 2419.3572
$ws.Cells.Item(132, 11).Value = 7258.071599999999
$ws.Cells.Item(132, 13).Value = -4728.071599999999

$ws.Cells.Item(136, 8).Value = 3873.75
$ws.Cells.Item(136, 9).Value = 3866.3572
$ws.Cells.Item(136, 11).Value = 11599.0716
$ws.Cells.Item(136, 13).Value = -9049.071599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(10, 8).Value = 10000
$ws.Cells.Item(10, 10).Value = 10000
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 14).Value = -10338

$ws.Cells.Item(96, 8).Value = 2966.5386
$ws.Cells.Item(96, 9).Value = 4148.125
$ws.Cells.Item(96, 10).Value = 1076
$ws.Cells.Item(96, 11).Value = 4148.125
$ws.Cells.Item(96, 12).Value = 1076
$ws.Cells.Item(96, 13).Value = -2775.125
$ws.Cells.Item(96, 14).Value = -3822

$ws.Cells.Item(132, 8).Value = 1373.6428
$ws.Cells.Item(132, 9).Value = 669.5
$ws.Cells.Item(132, 11).Value = 2008.5
$ws.Cells.Item(132, 13).Value = 521.5

$ws.Cells.Item(136, 8).Value = 1491.0968
$ws.Cells.Item(136, 9).Value = 1491.0968
$ws.Cells.Item(136, 11).Value = 4473.2904
$ws.Cells.Item(136, 13).Value = -1923.2904
